$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Helper: write a value into a cell as plain TEXT even when the
# string looks like a number (e.g. phone numbers / numeric-looking
# passwords), without leaving behind any unused number-format style.
# We do this by entering a text-literal formula ("=""...""") and then
# converting it in place to a static value via Copy / PasteSpecial
# (paste values only). This keeps the cell's style at the sheet
# default (no "s" attribute) and records the value as a shared string.
# -----------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
    $excel.CutCopyMode = $false
}

# ===================================================================
# Admin sheet: add the admin login row (email + password) and turn
# the email into a mailto hyperlink.
# ===================================================================
$admin = $wb.Worksheets.Item(1)
$admin.Cells.Item(2, 1).Value = "paramjotsingh966@gmail.com"
$admin.Cells.Item(2, 2).Value = 1234
$admin.Hyperlinks.Add($admin.Range("A2"), "mailto:paramjotsingh966@gmail.com")
$admin.Range("B3").Select()

# ===================================================================
# Users sheet: add the new registered user row.
# Columns: UserId | FirstName | LastName | Mobile | Email | Password | Interests
# ===================================================================
$users = $wb.Worksheets.Item(2)
$users.Activate()
$users.Cells.Item(2, 1).Value = "U#00001"
$users.Cells.Item(2, 2).Value = "Paramjot"
$users.Cells.Item(2, 3).Value = "Singh"
Set-TextValue $users.Range("D2") "9031398069"
$users.Cells.Item(2, 5).Value = "paramjotsingh966@gmail.com"
Set-TextValue $users.Range("F2") "1234"
$users.Cells.Item(2, 7).Value = "java,reactJs,networking,android"

# Column widths widen slightly to fit the new Mobile/Email content.
$users.Columns.Item(4).ColumnWidth = 33.666666666666664
$users.Columns.Item(5).ColumnWidth = 31.666666666666668

$users.Range("C3").Select()

# ===================================================================
# Resume sheet: no data change, just where the cursor ended up.
# ===================================================================
$resume = $wb.Worksheets.Item(5)
$resume.Activate()
$resume.Range("A10").Select()

# Users ends up the active tab.
$users.Activate()
